$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 13, pushing the existing
# rows 13..41 down to 14..42 (all other columns are identical across rows,
# only the date/volume/price columns differ row to row).
$ws.Rows("13:13").Insert()

$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 45044
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 100112043
$ws.Cells.Item(13, 7).Value = "Pepino dulce"
$ws.Cells.Item(13, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 15500
$ws.Cells.Item(13, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 861
$ws.Cells.Item(13, 17).Value = 18
$ws.Cells.Item(13, 18).Value = "Hortaliza"
